# pie chart and stacked bar chart for Metro view
# Populate the "Type" column (C) on the Modes sheet for the Bus and
# Demand Response modes so the upcoming pie/stacked-bar charts have a
# grouping column to work from.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Modes")

$ws.Range("C7").Value = "Bus"
$ws.Range("C8").Value = "Bus"
$ws.Range("C9").Value = "Bus"
$ws.Range("C12").Value = "Rail"
$ws.Range("C13").Value = "Demand"
$ws.Range("C15").Value = "Bus"
$ws.Range("C16").Value = "Demand"

$ws.Range("C13").Select()
